$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'25.654.73"
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.16%  '
$ws.Range('D3').Value = "'1.739.89"
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.44%  '
$ws.Range('D4').Value = "'1.004"
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').Value = "'237.88"
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -8.90%  '
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').Value = "'0.5022"
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.88%  '
$ws.Range('D8').Value = "'41.69"
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -6.73%  '
$ws.Range('D9').Value = "'0.2623"
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -13.13%  '
$ws.Range('D10').Value = "'0.06111"
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -11.15%  '
$ws.Range('D11').Value = "'1.748.75"
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.37%  '
$ws.Range('D12').Value = "'0.06952"
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -6.92%  '
$ws.Range('D13').Value = "'15.06"
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -15.10%  '
$ws.Range('D14').Value = "'4.464"
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -10.30%  '
$ws.Range('D15').Value = "'0.5927"
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -19.38%  '
$ws.Range('D16').Value = "'76.31"
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -14.73%  '
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('E18').Value = '  +0.17%  '
$ws.Range('D19').Value = "'25.712.22"
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.07%  '
$ws.Range('D20').Value = "'11.55"
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -17.21%  '
$ws.Range('D21').Value = "'0.000006752"
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -14.72%  '
$ws.Range('D22').Value = "'1.971.02"
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.47%  '
$ws.Range('D23').Value = "'4.027"
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -12.32%  '
$ws.Range('D24').Value = "'8.090"
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -12.86%  '
$ws.Range('D25').Value = "'5.067"
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -15.24%  '
$ws.Range('D26').Value = "'137.58"
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.73%  '
$ws.Range('D27').Value = "'1.532"
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -9.22%  '
$ws.Range('D28').Value = "'1.811"
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -18.40%  '
$ws.Range('D29').Value = "'14.93"
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -11.91%  '
$ws.Range('D30').Value = "'103.20"
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.69%  '
$ws.Range('D31').Value = "'3.750"
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -11.93%  '
$ws.Range('D32').Value = "'0.08085"
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -8.06%  '
$ws.Range('E33').Value = '  -14.84%  '
$ws.Range('D34').Value = "'0.04486"
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.37%  '
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('D36').Value = "'2.651"
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -9.18%  '
$ws.Range('D37').Value = "'0.9687"
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -14.26%  '
$ws.Range('D38').Value = "'0.6039"
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -17.50%  '
$ws.Range('D39').Value = "'2.645"
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -14.65%  '
$ws.Range('D40').Value = "'0.01541"
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -10.10%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = "'1.911"
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -16.56%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = "'1.002"
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = "'103.57"
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.79%  '
$ws.Range('D44').Value = "'5.118"
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -12.80%  '
$ws.Range('D45').Value = "'0.3777"
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -20.21%  '
$ws.Range('D46').Value = "'0.7255"
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -19.72%  '
$ws.Range('D47').Value = "'0.05329"
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -8.12%  '
$ws.Range('D48').Value = "'0.1101"
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -10.80%  '
$ws.Range('D49').Value = "'29.86"
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -14.49%  '
$ws.Range('D50').Value = "'5.842"
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -21.19%  '
$ws.Range('D51').Value = "'52.17"
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -13.32%  '
